{"js": "// Remove the \"Author\" paragraph (\"Ben Jarman\") entirely, per the diff:\n// the whole <w:p> with pStyle=\"Author\" is deleted (including its paragraph mark).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.style === \"Author\" && p.text.trim() === \"Ben Jarman\") {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Author\" paragraph (\"Ben Jarman\") entirely, per the diff:\n# the whole paragraph with style \"Author\" is deleted (including its mark).\n\n$d = $word.ActiveDocument\n\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Author\" -and $p.Range.Text.Trim() -eq \"Ben Jarman\") {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    $p.Range.Delete()\n}\n"}
